$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for numeric-looking strings so Excel
# does not auto-convert them to numbers (matches source data which
# stores prices/links/labels as plain text).

$ws.Range("D2").Value = '28.082.32'
$ws.Range("E2").Value = '  +1.54%  '

$ws.Range("D3").Value = '1.891.13'
$ws.Range("E3").Value = '  +1.17%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.013'
$ws.Range("E4").Value = '  +0.89%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '336.69'
$ws.Range("E5").Value = '  +1.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.012'
$ws.Range("E6").Value = '  +0.88%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4764'
$ws.Range("E7").Value = '  +1.63%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3966'
$ws.Range("E8").Value = '  +0.69%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.15'
$ws.Range("E9").Value = '  -1.49%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08053'
$ws.Range("E10").Value = '  -0.10%  '

$ws.Range("E11").Value = '  +0.08%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.00'
$ws.Range("E12").Value = '  +1.05%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.890.65'
$ws.Range("E13").Value = '  +0.90%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.047'
$ws.Range("E14").Value = '  +1.89%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.249'
$ws.Range("E15").Value = '  +1.56%  '

$ws.Range("E16").Value = '  +1.21%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.64'
$ws.Range("E17").Value = '  +2.30%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06768'
$ws.Range("E18").Value = '  +1.95%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001053'
$ws.Range("E19").Value = '  +0.61%  '

$ws.Range("E20").Value = '  +0.05%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.012'
$ws.Range("E21").Value = '  +0.94%  '

$ws.Range("D22").Value = '28.064.52'
$ws.Range("E22").Value = '  +1.40%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.540'
$ws.Range("E23").Value = '  +0.95%  '

$ws.Range("E24").Value = '  +0.55%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.354'
$ws.Range("E25").Value = '  +1.93%  '

$ws.Range("D26").Value = '2.104.37'
$ws.Range("E26").Value = '  +0.41%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.78'
$ws.Range("E27").Value = '  +1.33%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.05'
$ws.Range("E28").Value = '  -0.69%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.115'
$ws.Range("E29").Value = '  +1.22%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.554'
$ws.Range("E30").Value = '  +0.11%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '122.23'
$ws.Range("E31").Value = '  -0.15%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9815'
$ws.Range("E32").Value = '  +1.71%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09615'
$ws.Range("E33").Value = '  +1.35%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.638'
$ws.Range("E34").Value = '  +1.23%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.375'

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.373'
$ws.Range("E36").Value = '  -4.95%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02261'
$ws.Range("E37").Value = '  +0.09%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06089'
$ws.Range("E38").Value = '  +0.10%  '

$ws.Range("E39").Value = '  -2.13%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.234'
$ws.Range("E40").Value = '  +1.35%  '

$ws.Range("E41").Value = '  +0.91%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5989'
$ws.Range("E42").Value = '  +0.13%  '

$ws.Range("E43").Value = '  +0.39%  '

$ws.Range("E44").Value = '  +1.67%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.263'
$ws.Range("E45").Value = '  +0.91%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5679'
$ws.Range("E46").Value = '  -0.47%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.27'
$ws.Range("E47").Value = '  +0.45%  '

$ws.Range("E48").Value = '  +0.21%  '

$ws.Range("E49").Value = '  -0.66%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06830'
$ws.Range("E50").Value = '  -0.28%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '112.49'
$ws.Range("E51").Value = '  -1.54%  '
